$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
  2 = @{ 2=12.58994779371838; 3=8.040477012768791; 4=6.126848382458395; 6=30.12423707289883; 7=38.90214350870353; 8=16.81193118659591; 9=26.63419934908758; 11=9.889463519974987; 12=11.01558100805911; 13=14.842873163098; 14=20.73291257019853 }
  3 = @{ 2=12.37865971239469; 3=7.964693567312262; 4=6.100410957355625; 6=30.12565910354234; 7=38.89566614627927; 8=16.85065633576359; 9=26.70162062514474; 11=9.730037564699936; 12=11.02324636548798; 13=14.81590408389555; 14=20.79429618185272 }
  4 = @{ 2=12.24992547338738; 3=7.916949734598361; 4=6.083780789801942; 6=30.13368181655331; 7=38.90212889679466; 8=16.87707118599172; 9=26.74765762411512; 11=9.632966189882888; 12=11.02960154499391; 13=14.80179190116682; 14=20.83378575779217 }
  5 = @{ 2=12.19778760856527; 3=7.897195553627273; 4=6.076905231928023; 6=30.13874919876501; 7=38.90738682432981; 8=16.88849828565708; 9=26.76758351532482; 11=9.593666640250724; 12=11.03260638962504; 13=14.79666097424131; 14=20.85033185387068 }
  6 = @{ 2=12.18915172157894; 3=7.893897575583504; 4=6.075757632447275; 6=30.1396992189524; 7=38.90841833897402; 8=16.89043576860398; 9=26.77096254325024; 11=9.587158084342226; 12=11.03313042284706; 13=14.79584655409115; 14=20.85310676269359 }
  7 = @{ 2=12.24922092391282; 3=7.916684521186113; 4=6.083688461234097; 6=30.1337428772373; 7=38.90218918390836; 8=16.87722261260638; 9=26.74792163468419; 11=9.632435069209686; 12=11.02964038824107; 13=14.80172018798398; 14=20.83400706533486 }
  8 = @{ 2=12.51693145774518; 3=8.01460251737454; 4=6.117816042814496; 6=30.12324393889956; 7=38.89774332007374; 8=16.82473600124652; 9=26.65648237824481; 11=9.834355340658433; 12=11.01788217902374; 13=14.83306924778963; 14=20.75370481944388 }
  9 = @{ 2=13.04670945169527; 3=8.196712929464688; 4=6.181530022460087; 6=30.15933212682852; 7=38.97182329237975; 8=16.74275280960575; 9=26.51404477041187; 11=10.2345116285993; 12=11.00788055614694; 13=14.91374046291894; 14=20.61045705432221 }
  10 = @{ 2=13.43477067929405; 3=8.32410055119067; 4=6.226310855360278; 6=30.22026606098346; 7=39.07656236918482; 8=16.69530746725063; 9=26.43195089229962; 11=10.52804526375044; 12=11.0084539398717; 13=14.98439241893219; 14=20.51380509422977 }
  11 = @{ 2=13.61022568841254; 3=8.380585439687838; 4=6.246226783203035; 6=30.25540967851281; 7=39.13506070923553; 8=16.67650356909381; 9=26.39951616079822; 11=10.66086579809631; 12=11.01042467591382; 13=15.01892872427734; 14=20.4716843834547 }
  12 = @{ 2=13.6764452582505; 3=8.401757988096008; 4=6.253701649212987; 6=30.26977904406146; 7=39.15876334089234; 8=16.66978287134621; 9=26.38794104969791; 11=10.71101039955074; 12=11.01141574408444; 13=15.03234413921723; 14=20.45599869050415 }
  13 = @{ 2=13.6621945681453; 3=8.397207855653695; 4=6.252094799366893; 6=30.26663725895551; 7=39.15358977170987; 8=16.67121250298999; 9=26.39040248641865; 11=10.7002183758156; 12=11.01119143115848; 13=15.02944001047146; 14=20.45936513780327 }
  14 = @{ 2=13.61567842177145; 3=8.382331704505043; 4=6.246843093078244; 6=30.25657062666437; 7=39.13697972697061; 8=16.67594263664409; 9=26.3985496896478; 11=10.66499454105319; 12=11.01050131231947; 13=15.02002570261024; 14=20.47038861695344 }
  15 = @{ 2=13.58715517201474; 3=8.373191176551472; 4=6.243617521716702; 6=30.25054251950679; 7=39.12700719782664; 8=16.67889207016854; 9=26.40363221925712; 11=10.64339773701894; 12=11.01011043960038; 13=15.01430286490837; 14=20.47717523165633 }
  16 = @{ 2=13.42327686157052; 3=8.320379106049614; 4=6.22500002319402; 6=30.21811822203179; 7=39.07295693209283; 8=16.69659229670689; 9=26.43416948386858; 11=10.51934656526845; 12=11.00835946484589; 13=14.98218299377609; 14=20.5165948526484 }
  17 = @{ 2=13.32241886650038; 3=8.287600883568508; 4=6.213461215159289; 6=30.20012436586284; 7=39.04257187351251; 8=16.70816285725927; 9=26.45416157116084; 11=10.44302754721494; 12=11.0077225079341; 13=14.96308717326358; 14=20.54124975610944 }
  18 = @{ 2=13.26431055419739; 3=8.268610258964561; 4=6.206781646400637; 6=30.19047405634407; 7=39.02611735454605; 8=16.71507951833352; 9=26.46612260011933; 11=10.39906694522722; 12=11.00751716367173; 13=14.95232978492218; 14=20.55560449155482 }
  19 = @{ 2=13.24462141886859; 3=8.262156951512379; 4=6.204512760041286; 6=30.1873269055042; 7=39.02072196458295; 8=16.71746629816552; 9=26.47025173006033; 11=10.38417320751032; 12=11.00747532361862; 13=14.94872655452176; 14=20.56049465596058 }
  20 = @{ 2=13.33316595156756; 3=8.291104455190478; 4=6.21469397591911; 6=30.20196751228827; 7=39.0457006903913; 8=16.70690407648186; 9=26.45198554412796; 11=10.45115881859734; 12=11.00777365599542; 13=14.96509661531595; 14=20.53860721055165 }
  21 = @{ 2=13.62934786821257; 3=8.386707132172724; 4=6.248387472451916; 6=30.25949869538043; 7=39.14181650637074; 8=16.67454242578421; 9=26.39613745639878; 11=10.67534513644869; 12=11.01069738258987; 13=15.02278182010124; 14=20.46714358363271 }
  22 = @{ 2=13.82160037075646; 3=8.447920983140762; 4=6.270017736192424; 6=30.30328060629374; 7=39.21366571819411; 8=16.65572340520248; 9=26.36376015162291; 11=10.82095836493279; 12=11.01403446129637; 13=15.06244460604091; 14=20.42197923612924 }
  23 = @{ 2=13.71913359227417; 3=8.415368115527981; 4=6.258509444259548; 6=30.27935013084419; 7=39.17449580530423; 8=16.6655540891358; 9=26.38066296088373; 11=10.74334056003653; 12=11.01212328711178; 13=15.04109881388765; 14=20.4459436176495 }
  24 = @{ 2=13.32830757461037; 3=8.289520944718788; 4=6.214136787178429; 6=30.20113206166857; 7=39.04428299217707; 8=16.70747234708572; 9=26.45296787036412; 11=10.44748292795655; 12=11.00775003084564; 13=14.96418745745174; 14=20.53980134366089 }
  25 = @{ 2=12.90330568858274; 3=8.148546490868474; 4=6.164644907028731; 6=30.14351088703708; 7=38.9429277325049; 8=16.76268705648667; 9=26.54862207040489; 11=10.12612566644037; 12=11.00919197151198; 13=14.8898929772722; 14=20.64769479242103 }
}

foreach ($r in $data.Keys) {
  $rowData = $data[$r]
  foreach ($c in $rowData.Keys) {
    $ws.Cells.Item($r, $c).Value = $rowData[$c]
  }
}
